$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove "Split out Green Vegtables?" (was E4)
$ws.Range("E4").ClearContents()

# Remove "Move CCS to another File and tag elements" (was E12)
$ws.Range("E12").ClearContents()

# Remove "Timing chart printing" (was E19)
$ws.Range("E19").ClearContents()

# Add new entry "Clear Button" in E25
$ws.Range("E25").Value = "Clear Button"

# Reflect the new selection left in the saved file
[void]$ws.Range("E5:F6").Select()
